$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B, C, D hold text (date/time/weekday/week strings) in this sheet.
# Excel's COM layer auto-converts strings that "look like" a date/number
# (e.g. "2024-01-08" or "01") into real date/number values with an applied
# NumberFormat style. Force those two columns to be written as literal text
# by temporarily marking the cell as Text format, then clearing the format
# again afterwards so the new row matches the unstyled look of the other
# data rows (no numberFormat/style stamped on it).
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "2024-01-08"
$ws.Range("A32").ClearFormats()

$ws.Range("B32").Value = "10:31:46"

$ws.Range("C32").Value = "Monday"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "01"
$ws.Range("D32").ClearFormats()

$ws.Range("E32").Value = 139857
$ws.Range("F32").Value = 143220
$ws.Range("G32").Value = 171959
$ws.Range("H32").Value = 147215
$ws.Range("I32").Value = -1
$ws.Range("J32").Value = 118128
$ws.Range("K32").Value = 224704
$ws.Range("L32").Value = 249435
$ws.Range("M32").Value = 185402
$ws.Range("N32").Value = 110492
$ws.Range("O32").Value = 40648
$ws.Range("P32").Value = 30789
$ws.Range("Q32").Value = 72514
$ws.Range("R32").Value = -1
$ws.Range("S32").Value = 41977
$ws.Range("T32").Value = -1
